# Implements the LOQ4240.xlsx edit:
#  - Objetivos/Programa resumido/Programa/Método/Critério/Norma de recuperação/
#    Bibliografia rows get their real (previously-missing) Portuguese content,
#    shifting several existing values down by one row.
#  - A new "Bibliografia:" content row is appended at the bottom (row 22).
#  - Column A/B width definitions are split (col A alone vs. col B alone).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 22 needs the same look (styles) as row 21 before we put values in it ---
$ws.Range("A21:C21").Copy()
$ws.Range("A22:C22").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row 17 previously had no B/C cells; make sure they pick up the correct
#     column B/C styles (instead of defaulting to column A's style) ---
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C17").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Cell content updates (rows 10, 13-22) ---

# Row 10 - Objetivos: real text instead of the duplicated docente name
$ws.Range("B10").Value = "Introduzir os conceitos fundamentais da ciência administração e de configurações de uma organização."
$ws.Range("C10").Value = "Introduzir os conceitos fundamentais da ciência administração e de configurações de uma organização."

# Row 12 "Docentes responsáveis:" (A12) is unchanged.

# Row 13 - now holds only the docente name (B/C); A13 becomes blank
$ws.Range("A13").Clear()
$ws.Range("B13").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C13").Value = "11079086 - Herlandí de Souza Andrade"

# Row 14 - "Programa resumido:" with its real short-syllabus text
$ws.Range("A14").Value = "Programa resumido:"
$ws.Range("B14").Value = "1. Áreas de Atuação da Administração.2. Estrutura organizacional."
$ws.Range("C14").Value = "1. Áreas de Atuação da Administração.2. Estrutura organizacional."

# Row 15 - "Short syllabus:" (English) moves here
$ws.Range("A15").Value = "Short syllabus:"
$ws.Range("B15").Value = "1. Management Practice Areas. 2. Organizational structure"
$ws.Range("C15").Value = "1. Management Practice Areas. 2. Organizational structure"

# Row 16 - "Programa:" gets its real text
$ws.Range("A16").Value = "Programa:"
$ws.Range("B16").Value = "1. Noções básicas de Marketing, Finanças e Recursos Humanos. 2. Diferentes configurações de organização."
$ws.Range("C16").Value = "1. Noções básicas de Marketing, Finanças e Recursos Humanos. 2. Diferentes configurações de organização."

# Row 17 - "Syllabus:" (English) moves here
$ws.Range("A17").Value = "Syllabus:"
$ws.Range("B17").Value = "1. Basic notions of Marketing, Finance and Human Resources.2. Different organization settings."
$ws.Range("C17").Value = "1. Basic notions of Marketing, Finance and Human Resources.2. Different organization settings."

# Row 18 - "Avaliação:" header only, B/C cleared
$ws.Range("A18").Value = "Avaliação:"
$ws.Range("B18:C18").Clear()

# Row 19 - "Método:" with its real text (already present, now on this row)
$ws.Range("A19").Value = "Método:"
$ws.Range("B19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
$ws.Range("C19").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."

# Row 20 - "Critério:" with its real text
$ws.Range("A20").Value = "Critério:"
$ws.Range("B20").Value = "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas"
$ws.Range("C20").Value = "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas"

# Row 21 - "Norma de recuperação:" with its real formula text
$ws.Range("A21").Value = "Norma de recuperação:"
$ws.Range("B21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação"
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação"

# Row 22 (new) - "Bibliografia:" with the full bibliography text
$ws.Range("A22").Value = "Bibliografia:"
$ws.Range("B22").Value = "Chiavenato, I. Gestão de Pessoas. 4 ed. São Paulo: Manole, 2014.Chiavenato, I. Recursos Humanos: o capital humano das organizações. 10 ed. Rio de Janeiro, Campus, 2015.ROBBINS, S. P.; DECENZO, D. A.; WOLTER, R. Fundamentos de Gestão de Pessoas. São Paulo, saraiva, 2013.KOTLER, P. - ARMSTRONG, G. Princípios De Marketing. 15 ed. São Paulo: Pearson, 2014.KOTLER, P.; KELLER, K. L. Administração De Marketing. 15 ed. São Paulo: Pearson, 2019.CHIAVENATO, I. Introdução À Teoria Geral da Administração. 9 ed. São Paulo: Manole, 2014. MAXIMIANO, A. C. A. Teoria Geral da Administração: da Revolução Urbana À Revolução Digital. 8 ed. São Paulo: Atlas, 2017.GUERRINI, F. M.; ESCRIÇÃO FILHO, E.; ROSIM, D. Administração Para Engenheiros. Rio de Janeiro: Campus, 2016.CHIAVENATO, I. Administração Para Não Administradores: a Gestão de Negócios Ao Alcance de Todos. 2 ed. São Paulo: Manole, 2011.SILVA, M. M. L. Administração para Estudantes e Profissionais de Áreas Técnicas. São Paulo: Brasport, 2018.GITMAN, L. J. - ZUTTER, C. J. Princípios de Administração Financeira. 14 ed. São Paulo: Perason, 2017.GROPPELLI, A. A.; NIKBAKHT, E. Administração Financeira. 3 ed. São Paulo: Saraiva, 2010.MARCOUSÉ, I.; SURRIDGE, M.; GILLESPIE, A. Finanças. São Paulo: Saraiva, 2013.BOLMAN, L.G.; DEAL, T.E. Reframing organizations. San Francisco, John Wiley, 2013KOTLER, P.. O Marketing sem segredos. 1 ed. Porto Alegre. Bookman, 2005MINTZBERG, H. Criando organizações eficazes. 2 ed. São Paulo, Atlas, 2006.MORGAN, G. Imagens da organização. São Paulo, Atlas, 1996."
$ws.Range("C22").Value = "Chiavenato, I. Gestão de Pessoas. 4 ed. São Paulo: Manole, 2014.Chiavenato, I. Recursos Humanos: o capital humano das organizações. 10 ed. Rio de Janeiro, Campus, 2015.ROBBINS, S. P.; DECENZO, D. A.; WOLTER, R. Fundamentos de Gestão de Pessoas. São Paulo, saraiva, 2013.KOTLER, P. - ARMSTRONG, G. Princípios De Marketing. 15 ed. São Paulo: Pearson, 2014.KOTLER, P.; KELLER, K. L. Administração De Marketing. 15 ed. São Paulo: Pearson, 2019.CHIAVENATO, I. Introdução À Teoria Geral da Administração. 9 ed. São Paulo: Manole, 2014. MAXIMIANO, A. C. A. Teoria Geral da Administração: da Revolução Urbana À Revolução Digital. 8 ed. São Paulo: Atlas, 2017.GUERRINI, F. M.; ESCRIÇÃO FILHO, E.; ROSIM, D. Administração Para Engenheiros. Rio de Janeiro: Campus, 2016.CHIAVENATO, I. Administração Para Não Administradores: a Gestão de Negócios Ao Alcance de Todos. 2 ed. São Paulo: Manole, 2011.SILVA, M. M. L. Administração para Estudantes e Profissionais de Áreas Técnicas. São Paulo: Brasport, 2018.GITMAN, L. J. - ZUTTER, C. J. Princípios de Administração Financeira. 14 ed. São Paulo: Perason, 2017.GROPPELLI, A. A.; NIKBAKHT, E. Administração Financeira. 3 ed. São Paulo: Saraiva, 2010.MARCOUSÉ, I.; SURRIDGE, M.; GILLESPIE, A. Finanças. São Paulo: Saraiva, 2013.BOLMAN, L.G.; DEAL, T.E. Reframing organizations. San Francisco, John Wiley, 2013KOTLER, P.. O Marketing sem segredos. 1 ed. Porto Alegre. Bookman, 2005MINTZBERG, H. Criando organizações eficazes. 2 ed. São Paulo, Atlas, 2006.MORGAN, G. Imagens da organização. São Paulo, Atlas, 1996."

# --- Row height adjustments ---
$ws.Rows.Item(13).EntireRow.AutoFit()              # was 60 -> no explicit height now
$ws.Rows.Item(15).RowHeight = 60                    # was 120 -> 60
$ws.Rows.Item(17).RowHeight = 120                   # was none -> 120
$ws.Rows.Item(18).EntireRow.AutoFit()              # was 60 -> no explicit height now
$ws.Rows.Item(21).RowHeight = 60                    # was 120 -> 60
$ws.Rows.Item(22).RowHeight = 120                   # new row -> 120

# --- Column width/layout: split the merged A:B column-width range into
#     an A-only entry and keep B on its own (matches target structure) ---
$ws.Columns.Item(1).ColumnWidth = 29.83
$ws.Columns.Item(2).ColumnWidth = 59.83
